$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (some look numeric, e.g. "215.85"); keep them as
# text like the source data, not auto-converted numbers. Force text format,
# assign, then restore the default "Normal" style so no stray formatting diff
# is left behind on the cell.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.910.65'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.641.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5083'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.34%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2603'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06470'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.27'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.38%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07820'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.660.42'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.270'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.866.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.02%  '
$ws.Range('E15').Value = '  +2.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7717'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.56'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.918.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.004'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.78%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.402'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.58%  '
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.235'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.768'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '138.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.51%  '
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.874'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.65'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.246'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.05028'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.317'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.262'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.582'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.386'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9094'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.587'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5534'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.130.94'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01575'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9942'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.13'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.502'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('E44').Value = '  +1.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₈111'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4233'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.738'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05040'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.002'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.07%  '
